# Trade #9 closed at 2026-02-17 23:52:56 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet - roll up the new closed trade into the totals
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1500.26   # Current Capital
$summary.Range("B4").Value = 0.26      # Total P&L $
$summary.Range("B5").Value = 0.58      # Total P&L %
$summary.Range("B6").Value = 9         # Total Trades
$summary.Range("B7").Value = 5         # Winning Trades
$summary.Range("B9").Value = 55.56     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking strategy lives on row 6
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 100.26     # Capital
$status.Range("D6").Value = 9          # Trades
$status.Range("E6").Value = 0.26       # P&L $
$status.Range("F6").Value = 0.26       # P&L %
$status.Range("G6").Value = 55.56      # Win Rate %

# ---------------------------------------------------------------------
# Append the new trade (#9) to both the "All Trades" log and the
# per-strategy "MarketMaking" log - they mirror each other.
# ---------------------------------------------------------------------
$newRow = @(9, "2026-02-17", "23:52:50", "MarketMaking", "UP", 0.76, 0.89, "CLOSED", 17.1053, 0.13, 100.26, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.13)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 10

    # Column B holds a "YYYY-MM-DD" looking string that must stay plain
    # text (like the rest of the Date column) instead of being coerced
    # into a date serial by COM's type inference.
    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.NumberFormat = "@"

    for ($col = 1; $col -le $newRow.Length; $col++) {
        $ws.Cells.Item($row, $col).Value = $newRow[$col - 1]
    }

    # Drop the temporary text format now that the literal is safely
    # stored as a string, so the cell keeps the sheet's default style.
    $dateCell.Style = "Normal"
}
